$wb = $excel.ActiveWorkbook

$ssd = $wb.Worksheets.Item("SSD")
$usb = $wb.Worksheets.Item("USB")

$values = @(
    @(1416, 74),
    @(4993, 21),
    @(7489, 14),
    @(5518, 19),
    @(5242, 20),
    @(7489, 14),
    @(6553, 16),
    @(9532, 11),
    @(5242, 20),
    @(4369, 24)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $usb.Cells.Item($row, 2).Value = $values[$i][0]
    $usb.Cells.Item($row, 3).Value = $values[$i][1]
}

# Best-effort refresh of the USB line chart so its cached plot values (if the
# host supports recomputing chart caches) line up with the new data; harmless
# no-op if unsupported by the host.
try {
    $usbChart = $usb.ChartObjects(1).Chart
    $usbChart.SetSourceData($usb.Range("A1:B11"))
} catch {
}

$usb.Activate()
$usb.Range("B9").Select()
